$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 709.6957
$ws.Range("I19").Value = 366.6
$ws.Range("K19").Value = 366.6
$ws.Range("M19").Value = -191.6
$ws.Range("H51").Value = 2450.375
$ws.Range("I51").Value = 5750.5
$ws.Range("J51").Value = 1350.3334
$ws.Range("K51").Value = 5750.5
$ws.Range("L51").Value = 1350.3334
$ws.Range("M51").Value = -5266.5
$ws.Range("N51").Value = -2318.3334
$ws.Range("H129").Value = 859.4211
$ws.Range("I129").Value = 293.72726
$ws.Range("J129").Value = 955.1539
$ws.Range("K129").Value = 881.18178
$ws.Range("L129").Value = 2865.4617
$ws.Range("M129").Value = 4118.81822
$ws.Range("N129").Value = -12865.4617
$ws.Range("H138").Value = 3711.1804
$ws.Range("I138").Value = 2337.4736
$ws.Range("J138").Value = 4332.619
$ws.Range("K138").Value = 7012.4208
$ws.Range("L138").Value = 12997.857
$ws.Range("M138").Value = -1872.4208
$ws.Range("N138").Value = -23277.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16554.566
$ws.Range("I32").Value = 17582.82
$ws.Range("J32").Value = 9824.182000000001
$ws.Range("K32").Value = 17582.82
$ws.Range("L32").Value = 9824.182000000001
$ws.Range("M32").Value = -17295.82
$ws.Range("N32").Value = -10398.182
$ws.Range("H61").Value = 7552.096
$ws.Range("I61").Value = 5931.143
$ws.Range("J61").Value = 10889.353
$ws.Range("K61").Value = 5931.143
$ws.Range("L61").Value = 10889.353
$ws.Range("M61").Value = -5719.143
$ws.Range("N61").Value = -11313.353
$ws.Range("H63").Value = 4088.111
$ws.Range("I63").Value = 4129
$ws.Range("J63").Value = 3945
$ws.Range("K63").Value = 4129
$ws.Range("L63").Value = 3945
$ws.Range("M63").Value = -3443
$ws.Range("N63").Value = -5317
$ws.Range("H66").Value = 4088.111
$ws.Range("I66").Value = 4129
$ws.Range("J66").Value = 3945
$ws.Range("K66").Value = 20645
$ws.Range("L66").Value = 19725
$ws.Range("M66").Value = -17213
$ws.Range("N66").Value = -26589
$ws.Range("H88").Value = 13926.75
$ws.Range("I88").Value = 26000
$ws.Range("J88").Value = 1853.5
$ws.Range("K88").Value = 26000
$ws.Range("L88").Value = 1853.5
$ws.Range("M88").Value = -25594
$ws.Range("N88").Value = -2665.5
$ws.Range("H91").Value = 13926.75
$ws.Range("I91").Value = 26000
$ws.Range("J91").Value = 1853.5
$ws.Range("K91").Value = 26000
$ws.Range("L91").Value = 1853.5
$ws.Range("M91").Value = -24596
$ws.Range("N91").Value = -4661.5
$ws.Range("H122").Value = 5953723.5
$ws.Range("I122").Value = 1365.3529
$ws.Range("K122").Value = 4096.0587
$ws.Range("M122").Value = -1646.0587
$ws.Range("H132").Value = 4524.467
$ws.Range("I132").Value = 1761.2812
$ws.Range("J132").Value = 11326.154
$ws.Range("K132").Value = 5283.8436
$ws.Range("L132").Value = 33978.462
$ws.Range("M132").Value = -2753.8436
$ws.Range("N132").Value = -39038.462
$ws.Range("H136").Value = 7552.096
$ws.Range("I136").Value = 5931.143
$ws.Range("J136").Value = 10889.353
$ws.Range("K136").Value = 17793.429
$ws.Range("L136").Value = 32668.059
$ws.Range("M136").Value = -15243.429
$ws.Range("N136").Value = -37768.05899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1155.15
$ws.Range("I99").Value = 966.3214
$ws.Range("J99").Value = 1595.75
$ws.Range("K99").Value = 966.3214
$ws.Range("L99").Value = 1595.75
$ws.Range("M99").Value = 531.6786
$ws.Range("N99").Value = -4591.75
$ws.Range("H105").Value = 3937.8057
$ws.Range("I105").Value = 3342.92
$ws.Range("J105").Value = 5289.8184
$ws.Range("K105").Value = 3342.92
$ws.Range("L105").Value = 5289.8184
$ws.Range("M105").Value = -1595.92
$ws.Range("N105").Value = -8783.8184
$ws.Range("H134").Value = 41718.27
$ws.Range("I134").Value = 3474.1875
$ws.Range("K134").Value = 10422.5625
$ws.Range("M134").Value = -7887.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1453.4286
$ws.Range("I16").Value = 840
$ws.Range("J16").Value = 1913.5
$ws.Range("K16").Value = 840
$ws.Range("L16").Value = 1913.5
$ws.Range("M16").Value = -553
$ws.Range("N16").Value = -2487.5
$ws.Range("H31").Value = 5185.4
$ws.Range("I31").Value = 5291.4614
$ws.Range("J31").Value = 4879
$ws.Range("K31").Value = 5291.4614
$ws.Range("L31").Value = 4879
$ws.Range("M31").Value = -4996.4614
$ws.Range("N31").Value = -5469
$ws.Range("H34").Value = 5185.4
$ws.Range("I34").Value = 5291.4614
$ws.Range("J34").Value = 4879
$ws.Range("K34").Value = 5291.4614
$ws.Range("L34").Value = 4879
$ws.Range("M34").Value = -5089.4614
$ws.Range("N34").Value = -5283
$ws.Range("H107").Value = 770.4
$ws.Range("I107").Value = 946.1111
$ws.Range("J107").Value = 626.63635
$ws.Range("K107").Value = 946.1111
$ws.Range("L107").Value = 626.63635
$ws.Range("M107").Value = 973.8889
$ws.Range("N107").Value = -4466.63635
$ws.Range("H113").Value = 1453.4286
$ws.Range("I113").Value = 840
$ws.Range("J113").Value = 1913.5
$ws.Range("K113").Value = 840
$ws.Range("L113").Value = 1913.5
$ws.Range("M113").Value = 1330
$ws.Range("N113").Value = -6253.5
$ws.Range("H122").Value = 9782.666999999999
$ws.Range("I122").Value = 4706.5
$ws.Range("K122").Value = 14119.5
$ws.Range("M122").Value = -11669.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4389162.5
$ws.Range("I5").Value = 361.75
$ws.Range("J5").Value = 11912821
$ws.Range("K5").Value = 1085.25
$ws.Range("L5").Value = 35738463
$ws.Range("M5").Value = -973.25
$ws.Range("N5").Value = -35738687
$ws.Range("H98").Value = 537.8946999999999
$ws.Range("I98").Value = 423.07693
$ws.Range("J98").Value = 786.6667
$ws.Range("K98").Value = 1269.23079
$ws.Range("L98").Value = 2360.0001
$ws.Range("M98").Value = 228.7692099999999
$ws.Range("N98").Value = -5356.0001
$ws.Range("H122").Value = 732.1951
$ws.Range("I122").Value = 515.2
$ws.Range("K122").Value = 4636.8
$ws.Range("M122").Value = -2186.8
$ws.Range("H131").Value = 44892.91
$ws.Range("I131").Value = 1477.5
$ws.Range("J131").Value = 96991.39999999999
$ws.Range("K131").Value = 4432.5
$ws.Range("L131").Value = 290974.2
$ws.Range("M131").Value = 607.5
$ws.Range("N131").Value = -301054.2
$ws.Range("H135").Value = 4389162.5
$ws.Range("I135").Value = 361.75
$ws.Range("J135").Value = 11912821
$ws.Range("K135").Value = 3255.75
$ws.Range("L135").Value = 107215389
$ws.Range("M135").Value = -720.75
$ws.Range("N135").Value = -107220459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5495.25
$ws.Range("I102").Value = 5117.8423
$ws.Range("K102").Value = 5117.8423
$ws.Range("M102").Value = -3495.8423
$ws.Range("H122").Value = 6529.0835
$ws.Range("I122").Value = 7245
$ws.Range("K122").Value = 21735
$ws.Range("M122").Value = -19285
$ws.Range("H132").Value = 5696.1763
$ws.Range("I132").Value = 4268.718
$ws.Range("J132").Value = 10335.417
$ws.Range("K132").Value = 12806.154
$ws.Range("L132").Value = 31006.251
$ws.Range("M132").Value = -10276.154
$ws.Range("N132").Value = -36066.251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3464.4927
$ws.Range("I136").Value = 2116.9792
$ws.Range("K136").Value = 6350.937600000001
$ws.Range("M136").Value = -3800.937600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3003.4285
$ws.Range("I122").Value = 2428.7917
$ws.Range("K122").Value = 7286.375100000001
$ws.Range("M122").Value = -4836.375100000001
